$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")

# Pro number: 20581576 -> 20583536
$wsInput.Range("B3").Value = 20583536

# Instrument SN: A01606 -> A01612
$wsInput.Range("B4").Value = "A01612"

# ICA SN: (blank) -> APXCAS2134009
$wsInput.Range("B5").Value = "APXCAS2134009"
